$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (date format in col A, 2-decimal number format in col B)
# from the last existing data row (105) down into the two new rows (106-107)
$ws.Range("A105:B105").Copy() | Out-Null
$ws.Range("A106:B107").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# New data rows
$ws.Range("A106").Value = 45706
$ws.Range("B106").Value = 2.1

$ws.Range("A107").Value = 45716
$ws.Range("B107").Value = 1.86

# Update the view: scroll down a bit and move the active selection to A108
$ws.Range("A108").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 87
$excel.ActiveWindow.ScrollColumn = 1
